$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Checklist" to "Session"
$ws.Name = "Session"

# Keep the Student ID column stored as text (it already was text-as-number)
# so the new IDs don't silently become numeric cells.
$ws.Range("A2:A3").NumberFormat = "@"

# Update row 2: shift values from what was row 3, with new "Scan" type
$ws.Range("A2").Value = "231995"
$ws.Range("E2").Value = "Scan"

# Update row 3: shift values from what was row 4, with new "Scan" type
$ws.Range("A3").Value = "231996"
$ws.Range("D3").Value = "10:53:22"
$ws.Range("E3").Value = "Scan"

# Delete row 4 entirely, shrinking the used range to A1:F3
$ws.Rows.Item(4).Delete()
